$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1").Value = "Template_name_strategy"
$ws.Range("G1").Select()
